{"js": "// Add a new bulleted list item \"L\u00e0m silde b\u00e1o c\u00e1o\" right after the\n// existing list item that ends in \"... li\u00ean quan.\" (same list, numId 2)\n// and right before the following blank ListParagraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph to anchor on by its distinctive text rather than\n// a hard-coded index, so the script is resilient to minor doc changes.\nconst anchorText = \"T\u00ecm ki\u1ebfm c\u00e1c b\u00e0i b\u00e1o \u0111\u1ec1 t\u00e0i v\u1ec1 ph\u00e2n lo\u1ea1i d\u1eef li\u1ec7u\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for insertion.\");\n}\n\n// Insert the new list paragraph immediately after the anchor paragraph;\n// it inherits the surrounding list formatting (numId 2, ListParagraph\n// style) automatically because it is split off of that paragraph.\nconst newParagraph = anchor.insertParagraph(\"L\u00e0m silde b\u00e1o c\u00e1o\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new bulleted list item \"L\u00e0m silde b\u00e1o c\u00e1o\" right after the\n# existing list item that ends in \"... li\u00ean quan.\" (same list, numId 2)\n# and right before the following blank ListParagraph.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph to anchor on by its distinctive text rather than\n# a hard-coded index, so the script is resilient to minor doc changes.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*T\u00ecm ki\u1ebfm c\u00e1c b\u00e0i b\u00e1o \u0111\u1ec1 t\u00e0i v\u1ec1 ph\u00e2n lo\u1ea1i d\u1eef li\u1ec7u*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find anchor paragraph for insertion.\"\n}\n\n# Split a new paragraph mark after the anchor; the new (empty) paragraph\n# inherits the list/run formatting (numId 2, ListParagraph style) from\n# the anchor automatically.\n$target.Range.InsertParagraphAfter()\n\n# Fill in the text of the newly created paragraph.\n$newParagraph = $target.Next()\n$newParagraph.Range.InsertAfter(\"L\u00e0m silde b\u00e1o c\u00e1o\")\n"}
